# feat: add 2022-Q3 data
#
# The existing "2022-Q2" sheet becomes "2022-Q3" and gets the new quarter's
# fund-holding data; an exact copy of the old "2022-Q2" sheet (with its
# original data untouched) is inserted right after it, preserving the old
# quarter's numbers for the record. The "总计" (totals) sheet gets a new
# row for 2022-Q3, with the old 2022-Q2 totals row pushed down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate the current "2022-Q2" sheet so we keep an untouched archive
#    copy, then rename the two tabs into their final positions:
#      总计 | 2022-Q3 (was 2022-Q2, gets new data) | 2022-Q2 (archived copy)
# ---------------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($null, $wsQ2)
$wsArchive = $wb.Worksheets.Item("2022-Q2 (2)")

$wsQ2.Name = "2022-Q3"
$wsArchive.Name = "2022-Q2"

# Keep the first sheet active/selected, same as before the edit.
$wb.Worksheets.Item(1).Activate()

# ---------------------------------------------------------------------------
# 2. 总计 sheet: push the old 2022-Q2 totals row down to row 3 (bumping its
#    index column from 0 to 1), then write the new 2022-Q3 totals into row 2.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))
$wsTotal.Range("A3").Value = 1

$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 10
$wsTotal.Range("D2").Value = 0.67

# ---------------------------------------------------------------------------
# 3. 2022-Q3 sheet: replace the fund-holding table with the new quarter's
#    data entirely (new headers/rows, fresh formatting matching the "总计"
#    sheet's look rather than the old tab's).
# ---------------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")

# Wipe the old table (content + formatting) before rebuilding it.
$wsQ3.Range("A1:H11").Clear()

# Match "总计"'s page margins (0.75/0.75/1/1/0.5/0.5 in) rather than the old
# tab's (0.7/0.7/0.75/0.75/0.3/0.3 in).
$wsQ3.PageSetup.LeftMargin = 54
$wsQ3.PageSetup.RightMargin = 54
$wsQ3.PageSetup.TopMargin = 72
$wsQ3.PageSetup.BottomMargin = 72
$wsQ3.PageSetup.HeaderMargin = 36
$wsQ3.PageSetup.FooterMargin = 36

# Headers (B1:H1) and the index column (A2:A11) reuse "总计"'s bold/boxed
# header style — format-paint it over rather than re-deriving it, so no new
# style record gets minted.
$wsTotal.Range("A2").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)
$wsQ3.Range("A2:A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

# B (基金代码) and D:G (基金规模/股票总仓位/仓位占比/持有市值) are stored as
# text in the source data (keeps leading zeros / trailing zeros exactly as
# published), so force a text number format before writing them.
$wsQ3.Range("B2:B11").NumberFormat = "@"
$wsQ3.Range("D2:G11").NumberFormat = "@"

$q3rows = @(
    @(0, "001877", "宝盈国家安全沪港深股票A",   "8.74", "91.81", "3.33", "0.2910", 8),
    @(1, "013895", "宝盈成长精选混合A",         "8.51", "90.34", "2.94", "0.2502", 8),
    @(2, "013896", "宝盈成长精选混合C",         "2.93", "90.34", "2.94", "0.0861", 8),
    @(3, "501219", "华夏智胜先锋股票（LOF）A", "1.49", "92.63", "0.86", "0.0128", 6),
    @(4, "014198", "华夏智胜先锋股票（LOF）C", "1.26", "92.63", "0.86", "0.0108", 6),
    @(5, "013613", "宝盈国家安全沪港深股票C",   "0.23", "91.81", "3.33", "0.0077", 8),
    @(6, "013166", "东兴宸祥量化混合A",         "0.38", "93.87", "1.26", "0.0048", 4),
    @(7, "009327", "东兴兴晟混合A",             "0.38", "79.70", "1.08", "0.0041", 8),
    @(8, "013167", "东兴宸祥量化混合C",         "0.08", "93.87", "1.26", "0.0010", 4),
    @(9, "009328", "东兴兴晟混合C",             "0.07", "79.70", "1.08", "0.0008", 8)
)

$r = 2
foreach ($row in $q3rows) {
    $wsQ3.Range("A$r").Value = $row[0]
    $wsQ3.Range("B$r").Value = $row[1]
    $wsQ3.Range("C$r").Value = $row[2]
    $wsQ3.Range("D$r").Value = $row[3]
    $wsQ3.Range("E$r").Value = $row[4]
    $wsQ3.Range("F$r").Value = $row[5]
    $wsQ3.Range("G$r").Value = $row[6]
    $wsQ3.Range("H$r").Value = $row[7]
    $r = $r + 1
}
